$d = $word.ActiveDocument

# 1. Add a new run with "." after "It also acts as a container for all your elements"
$rng = $d.Content
$found = $rng.Find.Execute("It also acts as a container for all your elements", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(".")
    # Force the newly inserted text into its own run (matching formatting,
    # but a distinct <w:r>) by toggling a character property on/off.
    $rng.Bold = $true
    $rng.Bold = $false
}

# 2. Remove the stale <w:lastRenderedPageBreak/> marker preceding
#    "You can use opacity to give a tint to your headings etc." by
#    re-writing that run's text (Word drops the rendering-cache marker
#    whenever the run's text content is rewritten).
$rng2 = $d.Content
$rng2.Find.Execute("You can use opacity to give a tint to your headings etc.", $true, $true, $false, $false, $false, $true, 1, $false, "You can use opacity to give a tint to your headings etc.", 2)

# 3. Change "javascript" to "JSX"
$rng3 = $d.Content
$rng3.Find.Execute("javascript", $true, $true, $false, $false, $false, $true, 1, $false, "JSX", 2)
